$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing B10 value changes from "2" to "4.5" (kept as text, like the
# rest of column B, using a leading apostrophe so the numeric-looking
# value is stored as text rather than being parsed as a number)
$ws.Range("B10").Value = "'4.5"

# New row 11
$ws.Range("A11").Value = "'30/01/2024"
$ws.Range("B11").Value = "'3.5"

# New row 12
$ws.Range("A12").Value = "'31/01/2024"
$ws.Range("B12").Value = "'5"
